$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prime")

# Row 4: Max Agility
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "Valid Input"
$ws.Cells.Item(4,3).Value = "Max Agility"
$ws.Cells.Item(4,4).Value = 1
$ws.Cells.Item(4,5).Value = -2
$ws.Cells.Item(4,6).Value = 3
$ws.Cells.Item(4,7).Value = 1
$ws.Cells.Item(4,8).Value = 2
$ws.Cells.Item(4,9).Value = 3

# Row 5: Max Intelligence
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "Valid Input"
$ws.Cells.Item(5,3).Value = "Max Intelligence"
$ws.Cells.Item(5,4).Value = 1
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = -2
$ws.Cells.Item(5,7).Value = 3
$ws.Cells.Item(5,8).Value = 1
$ws.Cells.Item(5,9).Value = 3

# Row 6: Max Charisma
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "Valid Input"
$ws.Cells.Item(6,3).Value = "Max Charisma"
$ws.Cells.Item(6,4).Value = 1
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 2
$ws.Cells.Item(6,7).Value = -2
$ws.Cells.Item(6,8).Value = 3
$ws.Cells.Item(6,9).Value = 3

# Row 7: Even Attributes
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "Valid Input"
$ws.Cells.Item(7,3).Value = "Even Attributes"
$ws.Cells.Item(7,4).Value = 1
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 3
$ws.Cells.Item(7,7).Value = 0
$ws.Cells.Item(7,8).Value = -2
$ws.Cells.Item(7,9).Value = 3

$ws.Activate()
$ws.Range("I9").Select()
